$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Target cluster: ECs) ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5953293333333333
$ws.Range("H2").Value = 1.785988
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 32.56613966666666
$ws.Range("N2").Value = 97.698419
$ws.Range("O2").Value = 0.4203775945150413
$ws.Range("P2").Value = 0.4203775945150412
$ws.Range("Q2").Value = 19.38757821699689
$ws.Range("R2").Value = 174.488203952972
$ws.Range("S2").Value = 0.4203775945150413
$ws.Range("T2").Value = 0.4203775945150412

# --- Row 3 (Target cluster: FAPs) ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5953293333333333
$ws.Range("H3").Value = 1.785988
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.85529333333333
$ws.Range("N3").Value = 122.56588
$ws.Range("O3").Value = 0.5273775188114274
$ws.Range("P3").Value = 0.5273775188114271
$ws.Range("Q3").Value = 24.32235454327111
$ws.Range("R3").Value = 218.90119088944
$ws.Range("S3").Value = 0.5273775188114274
$ws.Range("T3").Value = 0.5273775188114271

# --- Row 4 (Target cluster: M1, was M2) ---
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5953293333333333
$ws.Range("H4").Value = 1.785988
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03779266666666666
$ws.Range("N4").Value = 0.113378
$ws.Range("O4").Value = 0.0004878438300104565
$ws.Range("P4").Value = 0.0004878438300104564
$ws.Range("Q4").Value = 0.02249908305155555
$ws.Range("R4").Value = 0.202491747464
$ws.Range("S4").Value = 0.0004878438300104565
$ws.Range("T4").Value = 0.0004878438300104564

# --- Row 5 (Target cluster: M2, was sCs) ---
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5953293333333333
$ws.Range("H5").Value = 1.785988
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.07240166666666666
$ws.Range("N5").Value = 0.217205
$ws.Range("O5").Value = 0.0009345915353721286
$ws.Range("P5").Value = 0.0009345915353721284
$ws.Range("Q5").Value = 0.04310283594888888
$ws.Range("R5").Value = 0.3879255235399999
$ws.Range("S5").Value = 0.0009345915353721286
$ws.Range("T5").Value = 0.0009345915353721284

# --- Row 6 (new row, Target cluster: sCs) ---
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Inha"
$ws.Range("C6").Value = "Tgfbr3"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5953293333333333
$ws.Range("H6").Value = 1.785988
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.937153333333333
$ws.Range("N6").Value = 11.81146
$ws.Range("O6").Value = 0.05082245130814891
$ws.Range("P6").Value = 0.05082245130814889
$ws.Range("Q6").Value = 2.343902869164444
$ws.Range("R6").Value = 21.09512582248
$ws.Range("S6").Value = 0.05082245130814891
$ws.Range("T6").Value = 0.05082245130814889
